$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 56: E56 becomes numeric, F56 timestamp updated slightly ---
$ws.Range("E56").Value = 919510038048
$ws.Range("F56").Value = 45965.83205196759

# --- Append new conversation rows 57-73 ---
# Row 57
$ws.Range("A57").Value = '[Call Started]'
$ws.Range("D57").Value = 'Good afternoon! I am your sales agent from Creer Infotech. Would you like to hear about our latest products?'
$ws.Range("E57").Value = 919510038048
$ws.Range("F57").Value = 45967.56912234954
$ws.Range("F57").NumberFormat = $ws.Range("F2").NumberFormat

# Row 58
$ws.Range("A58").Value = '[Call Started]'
$ws.Range("D58").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E58").Value = 919510038048
$ws.Range("F58").Value = 45972.81882967593
$ws.Range("F58").NumberFormat = $ws.Range("F2").NumberFormat

# Row 59
$ws.Range("A59").Value = '[Intro response]'
$ws.Range("B59").Value = 'Yes, I am still  there.'
$ws.Range("C59").Value = 'neutral'
$ws.Range("D59").Value = 'Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?'
$ws.Range("E59").Value = 919510038048
$ws.Range("F59").Value = 45972.81929960648
$ws.Range("F59").NumberFormat = $ws.Range("F2").NumberFormat

# Row 60
$ws.Range("A60").Value = '[Fallback]'
$ws.Range("B60").Value = 'Yes Bank'
$ws.Range("C60").Value = 'neutral'
$ws.Range("D60").Value = 'Sorry, we don’t have that product right now.
Here are our latest offers:
- Laptop Pro : one of the best laptop you can get right now with high end specs at ₹75000
- Smart watch : Best watch in market with all your daily tracking at ₹12000
- Bluetooth Earbuds: best anc earbuds with this price point  at ₹4000
Which product would you like to purchase?'
$ws.Range("E60").Value = 919510038048
$ws.Range("F60").Value = 45972.819420625
$ws.Range("F60").NumberFormat = $ws.Range("F2").NumberFormat

# Row 61
$ws.Range("A61").Value = '[Persuasion check]'
$ws.Range("B61").Value = 'none of the above no'
$ws.Range("C61").Value = 'neutral'
$ws.Range("D61").Value = 'I completely understand! But before you go — we’re giving a 20% discount just for today. Would you like to take a quick look?'
$ws.Range("E61").Value = 919510038048
$ws.Range("F61").Value = 45972.81985011574
$ws.Range("F61").NumberFormat = $ws.Range("F2").NumberFormat

# Row 62
$ws.Range("A62").Value = '[Fallback]'
$ws.Range("B62").Value = 'Hai yes, I am there right here.  Yes, I am there right  here.'
$ws.Range("C62").Value = 'neutral'
$ws.Range("D62").Value = 'Sorry, we don’t have that product right now.
Here are our latest offers:
- Laptop Pro : one of the best laptop you can get right now with high end specs at ₹75000
- Smart watch : Best watch in market with all your daily tracking at ₹12000
- Bluetooth Earbuds: best anc earbuds with this price point  at ₹4000
Which product would you like to purchase?'
$ws.Range("E62").Value = 919510038048
$ws.Range("F62").Value = 45972.82033612268
$ws.Range("F62").NumberFormat = $ws.Range("F2").NumberFormat

# Row 63
$ws.Range("A63").Value = '[Product match]'
$ws.Range("B63").Value = 'Yas Bluetooth earbuds?'
$ws.Range("C63").Value = 'neutral'
$ws.Range("D63").Value = 'Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 8048. Thank you for your time! I really appreciate it.'
$ws.Range("E63").Value = 919510038048
$ws.Range("F63").Value = 45972.82086135416
$ws.Range("F63").NumberFormat = $ws.Range("F2").NumberFormat

# Row 64
$ws.Range("A64").Value = '[Call Started]'
$ws.Range("D64").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E64").Value = 919510038048
$ws.Range("F64").Value = 45975.55740662037
$ws.Range("F64").NumberFormat = $ws.Range("F2").NumberFormat

# Row 65
$ws.Range("A65").Value = '[Intro response]'
$ws.Range("B65").Value = 'I M stil deyar. Hello I am speaking.'
$ws.Range("C65").Value = 'neutral'
$ws.Range("D65").Value = 'Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?'
$ws.Range("E65").Value = 919510038048
$ws.Range("F65").Value = 45975.55791690973
$ws.Range("F65").NumberFormat = $ws.Range("F2").NumberFormat

# Row 66
$ws.Range("A66").Value = '[Product match]'
$ws.Range("B66").Value = 'Bluetooth earbuds'
$ws.Range("C66").Value = 'neutral'
$ws.Range("D66").Value = 'Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 8048. Thank you for your time! I really appreciate it.'
$ws.Range("E66").Value = 919510038048
$ws.Range("F66").Value = 45975.55824385417
$ws.Range("F66").NumberFormat = $ws.Range("F2").NumberFormat

# Row 67
$ws.Range("A67").Value = '[Call Started]'
$ws.Range("D67").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E67").Value = 917990747606
$ws.Range("F67").Value = 45986.76252136574
$ws.Range("F67").NumberFormat = $ws.Range("F2").NumberFormat

# Row 68
$ws.Range("A68").Value = '[Persuasion check]'
$ws.Range("B68").Value = 'no no no'
$ws.Range("C68").Value = 'neutral'
$ws.Range("D68").Value = 'I completely understand! But before you go — we’re giving a 20% discount just for today. Would you like to take a quick look?'
$ws.Range("E68").Value = 917990747606
$ws.Range("F68").Value = 45986.76287478009
$ws.Range("F68").NumberFormat = $ws.Range("F2").NumberFormat

# Row 69
$ws.Range("A69").Value = '[Intro response]'
$ws.Range("B69").Value = 'yas yas'
$ws.Range("C69").Value = 'neutral'
$ws.Range("D69").Value = 'Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?'
$ws.Range("E69").Value = 917990747606
$ws.Range("F69").Value = 45986.76332363426
$ws.Range("F69").NumberFormat = $ws.Range("F2").NumberFormat

# Row 70
$ws.Range("A70").Value = '[Product match]'
$ws.Range("B70").Value = 'laptop Pro'
$ws.Range("C70").Value = 'neutral'
$ws.Range("D70").Value = 'Great choice! I’ve sent the link of Laptop Pro  to your phone number ending with 7606. Thank you for your time! I really appreciate it.'
$ws.Range("E70").Value = 917990747606
$ws.Range("F70").Value = 45986.76354986111
$ws.Range("F70").NumberFormat = $ws.Range("F2").NumberFormat

# Row 71
$ws.Range("A71").Value = '[Call Started]'
$ws.Range("D71").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E71").Value = 919510038048
$ws.Range("F71").Value = 45987.63015055555
$ws.Range("F71").NumberFormat = $ws.Range("F2").NumberFormat

# Row 72
$ws.Range("A72").Value = '[Intro response]'
$ws.Range("B72").Value = 'Yes.'
$ws.Range("C72").Value = 'neutral'
$ws.Range("D72").Value = 'Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?'
$ws.Range("E72").Value = 919510038048
$ws.Range("F72").Value = 45987.63036466435
$ws.Range("F72").NumberFormat = $ws.Range("F2").NumberFormat

# Row 73
$ws.Range("A73").Value = '[Product match]'
$ws.Range("B73").Value = 'Bluetooth earbuds'
$ws.Range("C73").Value = 'neutral'
$ws.Range("D73").Value = 'Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 8048. Thank you for your time! I really appreciate it.'
$ws.Range("E73").Value = "'" + '+919510038048'
$ws.Range("F73").Value = 45987.63049465394
$ws.Range("F73").NumberFormat = $ws.Range("F2").NumberFormat
